# corrrection typo fichiers de commentaires
# The column header "Décision du GT CNIG Paysage" contained a leftover
# typo from a copy/paste of another CNIG standard sheet ("Paysage"); this
# is a "risques" (ppr-risques) comments table, so the trailing word must
# read "Risques" instead.

$d = $word.ActiveDocument

$d.Content.Find.Execute("Paysage", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Risques", 2)
